$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) After ": UrbanSound8k" add: ",ESC50", "(AST già pronto ", the arxiv URL,
#    and the closing ")" as four separate runs.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$f1 = $r1.Find
$f1.Execute(": UrbanSound8k", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f1.Found) {
    throw "Could not find ': UrbanSound8k'"
}
$r1.Collapse(0)
$r1.InsertAfter(",ESC50")
$r1.Collapse(0)
$r1.InsertAfter("(AST già pronto ")
$r1.Collapse(0)
$r1.InsertAfter("https://arxiv.org/pdf/2104.01778.pdf")
$r1.Collapse(0)
$r1.InsertAfter(")")

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the "CNN14" run up to the
#    "Proprietà modelli " run (repagination shifted it by one run because of
#    the text inserted above). The marker has no textual representation in
#    the object model, so each paragraph is replaced wholesale (via
#    Range.InsertXML, which for this host replaces the exact range it is
#    called on) with the same paragraph content plus/minus the marker.
#    Document.Paragraphs is walked directly because a Range returned from a
#    Find on a sub-range doesn't reliably expose its own Paragraphs.
# ---------------------------------------------------------------------------

# --- 2a: add the marker to "Proprietà modelli " -----------------------------
$r2 = $d.Content
$f2 = $r2.Find
$f2.Execute("Proprietà modelli", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f2.Found) {
    throw "Could not find 'Proprietà modelli'"
}

$paras = $d.Paragraphs
$titlePara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $pr = $paras.Item($i).Range
    if ($pr.Start -le $r2.Start -and $pr.End -gt $r2.Start) {
        $titlePara = $pr
    }
}
if ($null -eq $titlePara) {
    throw "Could not resolve the 'Proprietà modelli' paragraph"
}

$titleXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2FF5A65A" w14:textId="77777777" w:rsidR="00A403DC" w:rsidRPr="00FF3839" w:rsidRDefault="00A403DC" w:rsidP="00A403DC"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Proprietà modelli </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.InsertXML($titleXml)

# --- 2b: remove the marker from the "CNN14" run ------------------------------
$r3 = $d.Content
$f3 = $r3.Find
$f3.Execute("Proprietà modelli", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f3.Found) {
    throw "Could not find 'Proprietà modelli' (second pass)"
}
$afterTitle = $d.Range($r3.End, $d.Content.End)
$f4 = $afterTitle.Find
$f4.Execute("CNN14", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f4.Found) {
    throw "Could not find 'CNN14' after the title"
}

$cnnPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $pr = $paras.Item($i).Range
    if ($pr.Start -le $afterTitle.Start -and $pr.End -gt $afterTitle.Start) {
        $cnnPara = $pr
    }
}
if ($null -eq $cnnPara) {
    throw "Could not resolve the 'CNN14' paragraph"
}

$cnnXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7DFBA2CF" w14:textId="77777777" w:rsidR="00A403DC" w:rsidRDefault="00A403DC" w:rsidP="00A403DC"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00C36C03"><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>CNN14</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:Rete convoluzionale creata per la classificazione audio.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cnnPara.InsertXML($cnnXml)
